$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4006370.2
$ws.Range("J17").Value = 4006370.2
$ws.Range("L17").Value = 12019110.6
$ws.Range("N17").Value = -12019446.6
$ws.Range("H40").Value = 974.76
$ws.Range("I40").Value = 704.05884
$ws.Range("J40").Value = 1550
$ws.Range("K40").Value = 704.05884
$ws.Range("L40").Value = 1550
$ws.Range("M40").Value = -529.05884
$ws.Range("N40").Value = -1900
$ws.Range("H70").Value = 1580
$ws.Range("I70").Value = 1602.5
$ws.Range("J70").Value = 1490
$ws.Range("K70").Value = 4807.5
$ws.Range("L70").Value = 4470
$ws.Range("M70").Value = -4537.5
$ws.Range("N70").Value = -5010
$ws.Range("H73").Value = 1580
$ws.Range("I73").Value = 1602.5
$ws.Range("J73").Value = 1490
$ws.Range("K73").Value = 4807.5
$ws.Range("L73").Value = 4470
$ws.Range("M73").Value = -3871.5
$ws.Range("N73").Value = -6342
$ws.Range("H86").Value = 17653.5
$ws.Range("I86").Value = 931.6667
$ws.Range("K86").Value = 931.6667
$ws.Range("M86").Value = 191.3333
$ws.Range("H89").Value = 17653.5
$ws.Range("I89").Value = 931.6667
$ws.Range("K89").Value = 4658.3335
$ws.Range("M89").Value = 957.6665000000003
$ws.Range("H132").Value = 17545658
$ws.Range("I132").Value = 18520398
$ws.Range("J132").Value = 325.33334
$ws.Range("K132").Value = 55561194
$ws.Range("L132").Value = 976.0000200000001
$ws.Range("M132").Value = -55558664
$ws.Range("N132").Value = -6036.00002
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 55571028
$ws.Range("I135").Value = 2639.8
$ws.Range("K135").Value = 23758.2
$ws.Range("M135").Value = -21223.2
$ws.Range("H137").Value = 120197.12
$ws.Range("I137").Value = 140062.83
$ws.Range("K137").Value = 420188.49
$ws.Range("M137").Value = -417638.49
$ws.Range("H138").Value = 4517.6836
$ws.Range("I138").Value = 3678.3572
$ws.Range("J138").Value = 4657.5713
$ws.Range("K138").Value = 11035.0716
$ws.Range("L138").Value = 13972.7139
$ws.Range("M138").Value = -5895.071599999999
$ws.Range("N138").Value = -24252.7139
$ws.Range("H141").Value = 1483.2069
$ws.Range("I141").Value = 1537.5186
$ws.Range("J141").Value = 750
$ws.Range("K141").Value = 4612.5558
$ws.Range("L141").Value = 2250
$ws.Range("M141").Value = 567.4441999999999
$ws.Range("N141").Value = -12610

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22923.924
$ws.Range("I32").Value = 18501.375
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 18501.375
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -18214.375
$ws.Range("N32").Value = -30574
$ws.Range("H132").Value = 7822762
$ws.Range("I132").Value = 9436243
$ws.Range("J132").Value = 48718.91
$ws.Range("K132").Value = 28308729
$ws.Range("L132").Value = 146156.73
$ws.Range("M132").Value = -28306199
$ws.Range("N132").Value = -151216.73

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3986.2683
$ws.Range("I134").Value = 3783.0293
$ws.Range("J134").Value = 4973.4287
$ws.Range("K134").Value = 11349.0879
$ws.Range("L134").Value = 14920.2861
$ws.Range("M134").Value = -8814.0879
$ws.Range("N134").Value = -19990.2861

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5984.3335
$ws.Range("I31").Value = 3772.389
$ws.Range("J31").Value = 8638.666999999999
$ws.Range("K31").Value = 3772.389
$ws.Range("L31").Value = 8638.666999999999
$ws.Range("M31").Value = -3477.389
$ws.Range("N31").Value = -9228.666999999999
$ws.Range("H34").Value = 5984.3335
$ws.Range("I34").Value = 3772.389
$ws.Range("J34").Value = 8638.666999999999
$ws.Range("K34").Value = 3772.389
$ws.Range("L34").Value = 8638.666999999999
$ws.Range("M34").Value = -3570.389
$ws.Range("N34").Value = -9042.666999999999
$ws.Range("H99").Value = 3792.2
$ws.Range("I99").Value = 2840.25
$ws.Range("K99").Value = 2840.25
$ws.Range("M99").Value = -1342.25
$ws.Range("H126").Value = 3792.2
$ws.Range("I126").Value = 2840.25
$ws.Range("K126").Value = 8520.75
$ws.Range("M126").Value = -6050.75
$ws.Range("H132").Value = 37040550
$ws.Range("I132").Value = 45457268
$ws.Range("J132").Value = 7002.4
$ws.Range("K132").Value = 136371804
$ws.Range("L132").Value = 21007.2
$ws.Range("M132").Value = -136369274
$ws.Range("N132").Value = -26067.2
$ws.Range("H134").Value = 30303734
$ws.Range("I134").Value = 32258748
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 96776244
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -96773709
$ws.Range("N134").Value = -8070
$ws.Range("H141").Value = 41198
$ws.Range("J141").Value = 41198
$ws.Range("L141").Value = 41198
$ws.Range("N141").Value = -51558

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1641.4736
$ws.Range("J122").Value = 1788.1177
$ws.Range("L122").Value = 16093.0593
$ws.Range("N122").Value = -20993.0593
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H131").Value = 757.76
$ws.Range("J131").Value = 784.54346
$ws.Range("L131").Value = 2353.63038
$ws.Range("N131").Value = -12433.63038
$ws.Range("H132").Value = 1981
$ws.Range("I132").Value = 1002.5
$ws.Range("J132").Value = 2633.3333
$ws.Range("K132").Value = 9022.5
$ws.Range("L132").Value = 23699.9997
$ws.Range("M132").Value = -6492.5
$ws.Range("N132").Value = -28759.9997
$ws.Range("H133").Value = 6772.3335
$ws.Range("I133").Value = 1960
$ws.Range("J133").Value = 7975.4165
$ws.Range("K133").Value = 5880
$ws.Range("L133").Value = 23926.2495
$ws.Range("M133").Value = -820
$ws.Range("N133").Value = -34046.24950000001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3800.1428
$ws.Range("I80").Value = 3474.75
$ws.Range("J80").Value = 4000.3845
$ws.Range("K80").Value = 3474.75
$ws.Range("L80").Value = 4000.3845
$ws.Range("M80").Value = -2476.75
$ws.Range("N80").Value = -5996.3845
$ws.Range("H83").Value = 3800.1428
$ws.Range("I83").Value = 3474.75
$ws.Range("J83").Value = 4000.3845
$ws.Range("K83").Value = 17373.75
$ws.Range("L83").Value = 20001.9225
$ws.Range("M83").Value = -12381.75
$ws.Range("N83").Value = -29985.9225

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5422.143
$ws.Range("I7").Value = 7157.143
$ws.Range("J7").Value = 3687.1428
$ws.Range("K7").Value = 7157.143
$ws.Range("L7").Value = 3687.1428
$ws.Range("M7").Value = -7045.143
$ws.Range("N7").Value = -3911.1428
$ws.Range("H16").Value = 551.9375
$ws.Range("I16").Value = 584.2143
$ws.Range("J16").Value = 326
$ws.Range("K16").Value = 584.2143
$ws.Range("L16").Value = 326
$ws.Range("M16").Value = -414.2143
$ws.Range("N16").Value = -666
$ws.Range("H68").Value = 3049.1667
$ws.Range("I68").Value = 2934
$ws.Range("K68").Value = 2934
$ws.Range("M68").Value = -2185
$ws.Range("H71").Value = 3049.1667
$ws.Range("I71").Value = 2934
$ws.Range("K71").Value = 14670
$ws.Range("M71").Value = -10926
$ws.Range("H126").Value = 5422.143
$ws.Range("I126").Value = 7157.143
$ws.Range("J126").Value = 3687.1428
$ws.Range("K126").Value = 21471.429
$ws.Range("L126").Value = 11061.4284
$ws.Range("M126").Value = -19001.429
$ws.Range("N126").Value = -16001.4284
$ws.Range("H132").Value = 1890.4333
$ws.Range("I132").Value = 1100.1666
$ws.Range("K132").Value = 3300.4998
$ws.Range("M132").Value = -770.4998000000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 8000
$ws.Range("J26").Value = 8000
$ws.Range("L26").Value = 8000
$ws.Range("N26").Value = -8586
$ws.Range("H62").Value = 4999.5
$ws.Range("I62").Value = 4999.5
$ws.Range("K62").Value = 4999.5
$ws.Range("M62").Value = -4375.5
$ws.Range("H65").Value = 4999.5
$ws.Range("I65").Value = 4999.5
$ws.Range("K65").Value = 24997.5
$ws.Range("M65").Value = -21877.5
$ws.Range("H126").Value = 2160.4546
$ws.Range("I126").Value = 1698.5714
$ws.Range("J126").Value = 2968.75
$ws.Range("K126").Value = 5095.7142
$ws.Range("L126").Value = 8906.25
$ws.Range("M126").Value = -2625.7142
$ws.Range("N126").Value = -13846.25
$ws.Range("H132").Value = 15626537
$ws.Range("I132").Value = 26316760
$ws.Range("J132").Value = 2362.5386
$ws.Range("K132").Value = 78950280
$ws.Range("L132").Value = 7087.6158
$ws.Range("M132").Value = -78947750
$ws.Range("N132").Value = -12147.6158
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 32261738
$ws.Range("I136").Value = 50001972
$ws.Range("K136").Value = 150005916
$ws.Range("M136").Value = -150003366
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 76139
$ws.Range("J141").Value = 76139
$ws.Range("L141").Value = 76139
$ws.Range("N141").Value = -86499
